$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.538.73'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -2.14%  '
$c.Style = "Normal"

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.580.04'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -2.98%  '
$c.Style = "Normal"

# Row 4
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.27%  '
$c.Style = "Normal"

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '210.44'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -2.66%  '
$c.Style = "Normal"

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.504'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -2.04%  '
$c.Style = "Normal"

# Row 7
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  +0.19%  '
$c.Style = "Normal"

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.247'
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -2.08%  '
$c.Style = "Normal"

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.0615'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -1.25%  '
$c.Style = "Normal"

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.38'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -3.69%  '
$c.Style = "Normal"

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0835'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -1.91%  '
$c.Style = "Normal"

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.809.75'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -2.45%  '
$c.Style = "Normal"

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.593.23'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -1.50%  '
$c.Style = "Normal"

# Row 14
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -1.44%  '
$c.Style = "Normal"

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.525'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -2.49%  '
$c.Style = "Normal"

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '63.64'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -1.55%  '
$c.Style = "Normal"

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '26.592.66'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -1.82%  '
$c.Style = "Normal"

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.0₃0726'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -0.67%  '
$c.Style = "Normal"

# Row 19
$c = $ws.Range("B19")
$c.NumberFormat = "@"
$c.Value = 'Dai'
$c.Style = "Normal"
$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +0.06%  '
$c.Style = "Normal"

# Row 20
$c = $ws.Range("B20")
$c.NumberFormat = "@"
$c.Value = 'BitcoinCash'
$c.Style = "Normal"
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '207.92'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -2.70%  '
$c.Style = "Normal"

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.65'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -3.32%  '
$c.Style = "Normal"

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.24'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -3.15%  '
$c.Style = "Normal"

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -5.54%  '
$c.Style = "Normal"

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '8.86'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -1.87%  '
$c.Style = "Normal"

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '146.19'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -1.40%  '
$c.Style = "Normal"

# Row 26
$c = $ws.Range("B26")
$c.NumberFormat = "@"
$c.Value = 'BinanceUSD'
$c.Style = "Normal"
$c = $ws.Range("C26")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +0.45%  '
$c.Style = "Normal"

# Row 27
$c = $ws.Range("B27")
$c.NumberFormat = "@"
$c.Value = 'Cosmos'
$c.Style = "Normal"
$c = $ws.Range("C27")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.44'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +2.18%  '
$c.Style = "Normal"

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.112'
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -4.55%  '
$c.Style = "Normal"

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '15.23'
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -1.75%  '
$c.Style = "Normal"

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0498'
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -1.10%  '
$c.Style = "Normal"

# Row 31
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -2.10%  '
$c.Style = "Normal"

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.23'
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -4.27%  '
$c.Style = "Normal"

# Row 33
$c = $ws.Range("B33")
$c.NumberFormat = "@"
$c.Value = 'ImmutableX'
$c.Style = "Normal"
$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.649'
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +21.63%  '
$c.Style = "Normal"

# Row 34
$c = $ws.Range("B34")
$c.NumberFormat = "@"
$c.Value = 'InternetComputer(DFINITY)'
$c.Style = "Normal"
$c = $ws.Range("C34")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.91'
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -3.07%  '
$c.Style = "Normal"

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.301.70'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -0.89%  '
$c.Style = "Normal"

# Row 36
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -3.69%  '
$c.Style = "Normal"

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.42'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -1.21%  '
$c.Style = "Normal"

# Row 38
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -1.43%  '
$c.Style = "Normal"

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.814'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -3.29%  '
$c.Style = "Normal"

# Row 40
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.Style = "Normal"

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.783'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -2.26%  '
$c.Style = "Normal"

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.16'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -4.11%  '
$c.Style = "Normal"

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.27'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  +0.65%  '
$c.Style = "Normal"

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '62.63'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -1.64%  '
$c.Style = "Normal"

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.723.69'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -2.25%  '
$c.Style = "Normal"

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '88.71'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -2.12%  '
$c.Style = "Normal"

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.60'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +0.29%  '
$c.Style = "Normal"

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.828'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +4.85%  '
$c.Style = "Normal"

# Row 49
$c = $ws.Range("B49")
$c.NumberFormat = "@"
$c.Value = 'BabyDogeCoin'
$c.Style = "Normal"
$c = $ws.Range("C49")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.0₆0103'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -3.16%  '
$c.Style = "Normal"

# Row 50
$c = $ws.Range("B50")
$c.NumberFormat = "@"
$c.Value = 'Cronos'
$c.Style = "Normal"
$c = $ws.Range("C50")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0506'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -1.45%  '
$c.Style = "Normal"

# Row 51
$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = 'Algorand'
$c.Style = "Normal"
$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0973'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +2.96%  '
$c.Style = "Normal"
